# Fix: corregir el endpoint de cargar el stock del productos en megabahia
# The source feed now yields a single product record (row 2) with updated
# data, and the two extra product rows (3 and 4) that used to be appended
# are no longer produced, so they are removed from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new product's data.
$ws.Cells.Item(2, 1).Value2 = "HOAC01191"
$ws.Cells.Item(2, 2).Value2 = "ESPUMA DE PUERTA DOBLE"
$ws.Cells.Item(2, 3).Value2 = "86"
$ws.Cells.Item(2, 4).Value2 = "86"
$ws.Cells.Item(2, 6).Value2 = "5.95"
$ws.Cells.Item(2, 7).Value2 = "5.75"
$ws.Cells.Item(2, 8).Value2 = "5.5"
$ws.Cells.Item(2, 9).Value2 = "4.95"
$ws.Cells.Item(2, 10).Value2 = "4.75"
$ws.Cells.Item(2, 13).Value2 = "ACCESORIOS PARA EL HOGAR"
$ws.Cells.Item(2, 29).Value2 = "HOAC"

# Remove the now-absent extra product rows (3 and 4) entirely.
$ws.Rows.Item(4).EntireRow.Delete() | Out-Null
$ws.Rows.Item(3).EntireRow.Delete() | Out-Null
